# The commit swaps the deck's active theme ("Integral" / Red Violet colour
# scheme, ppt/theme/theme1.xml) for the stock PowerPoint "Office Theme"
# colour scheme (previously only used by the notes master, ppt/theme/theme2.xml).
#
# Only the <a:clrScheme> colour values (and the theme/clrScheme display
# names) differ between the two theme parts - fontScheme/fmtScheme are
# byte-identical - so the visible effect of the edit is entirely captured
# by re-pointing the presentation's 12 theme colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) at the "Office" palette.
#
# PowerPoint's RGB() helper packs red/green/blue into a single BGR-ordered
# long, which is what ColorScheme/ThemeColorScheme.Item(n).RGB expects.
function Office-RGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Slide.ThemeColorScheme reaches into the same clrScheme element as
# SlideMaster.ColorScheme but (unlike the master's own ColorScheme
# collection) leaves the clrScheme's name attribute alone instead of
# blanking it out, so prefer it here.
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = Office-RGB 0x00 0x00 0x00   # dk1      #000000
$tcs.Item(2).RGB  = Office-RGB 0xFF 0xFF 0xFF   # lt1      #FFFFFF
$tcs.Item(3).RGB  = Office-RGB 0x44 0x54 0x6A   # dk2      #44546A
$tcs.Item(4).RGB  = Office-RGB 0xE7 0xE6 0xE6   # lt2      #E7E6E6
$tcs.Item(5).RGB  = Office-RGB 0x5B 0x9B 0xD5   # accent1  #5B9BD5
$tcs.Item(6).RGB  = Office-RGB 0xED 0x7D 0x31   # accent2  #ED7D31
$tcs.Item(7).RGB  = Office-RGB 0xA5 0xA5 0xA5   # accent3  #A5A5A5
$tcs.Item(8).RGB  = Office-RGB 0xFF 0xC0 0x00   # accent4  #FFC000
$tcs.Item(9).RGB  = Office-RGB 0x44 0x72 0xC4   # accent5  #4472C4
$tcs.Item(10).RGB = Office-RGB 0x70 0xAD 0x47   # accent6  #70AD47
$tcs.Item(11).RGB = Office-RGB 0x05 0x63 0xC1   # hlink    #0563C1
$tcs.Item(12).RGB = Office-RGB 0x95 0x4F 0x72   # folHlink #954F72
